# "Update view survey inprogress"
#
# The "my_tasks" row's texts change from a generic "My Tasks" label to a
# survey-specific call-to-action, and the "in_progress_surveys" row's texts
# change from a dashboard subtitle to a "survey history" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 -> variable "my_tasks"
$ws.Range("B10").Value = "Thực hiện khảo sát"
$ws.Range("C10").Value = "Take a pollution survey"

# Row 16 -> variable "in_progress_surveys"
$ws.Range("C16").Value = "Survey history"
$ws.Range("B16").Value = "Lịch sử các tiến trình khảo sát"

# Restore the cursor/selection the author left the sheet on.
$ws.Range("B17").Select()
